$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'schubert-winterreise_119'
$ws.Range("B2").Value = 'schubert-winterreise_200'
$ws.Range("C2").Value = 0.7083333333333333
$ws.Range("D2").Value = '[[''D:7'', ''G:maj'', ''D:7'', ''G:maj'', ''D:7'', ''G:maj'']]'
$ws.Range("E2").Value = '[[''B:7'', ''E:maj'', ''B:7'', ''E:maj'', ''B:7'', ''E:maj'']]'
$ws.Range("F2").Value = '[(34.3, 49.84)]'
$ws.Range("G2").Value = '[(16.34, 44.36)]'
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""

$ws.Range("A3").Value = 'schubert-winterreise_150'
$ws.Range("B3").Value = 'schubert-winterreise_98'
$ws.Range("C3").Value = 0.3736263736263736
$ws.Range("D3").Value = '[[''A:maj'', ''D:min'', ''A:maj'', ''D:min'', ''A:maj'', ''D:min'']]'
$ws.Range("E3").Value = '[[''C'', ''F:min/C'', ''C'', ''F:min'', ''C'', ''F:min'']]'
$ws.Range("F3").Value = '[(0.78, 9.74)]'
$ws.Range("G3").Value = '[(45.26, 57.58)]'
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""

$ws.Range("A4").Value = 'schubert-winterreise_203'
$ws.Range("B4").Value = 'isophonics_19'
$ws.Range("C4").Value = 0.1230769230769231
$ws.Range("D4").Value = '[[''G:maj'', ''G:7/F'', ''C:maj/E'']]'
$ws.Range("E4").Value = '[[''F'', ''F:7'', ''Bb'']]'
$ws.Range("F4").Value = '[(54.4, 61.44)]'
$ws.Range("G4").Value = '[(38.438956, 41.94517)]'
$ws.Range("H4").Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'
$ws.Range("I4").Value = ""

$ws.Range("A5").Value = 'jaah_14'
$ws.Range("B5").Value = 'schubert-winterreise_33'
$ws.Range("C5").Value = 0.2597402597402597
$ws.Range("D5").Value = '[[''Ab'', ''Eb:7'', ''Ab'']]'
$ws.Range("E5").Value = '[[''G:maj/D'', ''D:7'', ''G:maj'']]'
$ws.Range("F5").Value = '[(66.93, 72.56)]'
$ws.Range("G5").Value = '[(65.44, 66.82)]'
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""

$ws.Range("A6").Value = 'isophonics_156'
$ws.Range("B6").Value = 'jaah_79'
$ws.Range("C6").Value = 0.1672201138519924
$ws.Range("D6").Value = '[[''A'', ''A'', ''A'', ''A:7/3'', ''D'']]'
$ws.Range("E6").Value = '[[''Bb'', ''Bb'', ''Bb'', ''Bb:7'', ''Eb'']]'
$ws.Range("F6").Value = '[(14.54151, 21.658426)]'
$ws.Range("G6").Value = '[(11.95, 15.8)]'
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""

$ws.Range("A7").Value = 'schubert-winterreise_70'
$ws.Range("B7").Value = 'schubert-winterreise_171'
$ws.Range("C7").Value = 0.2015810276679842
$ws.Range("D7").Value = '[[''B:min'', ''F#:7/C#'', ''B:min/D'']]'
$ws.Range("E7").Value = '[[''F#:min'', ''C#:7'', ''F#:min'']]'
$ws.Range("F7").Value = '[(1.54, 3.58)]'
$ws.Range("G7").Value = '[(3.82, 6.46)]'
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

$ws.Range("A8").Value = 'isophonics_242'
$ws.Range("B8").Value = 'isophonics_200'
$ws.Range("C8").Value = 0.1169415292353823
$ws.Range("D8").Value = '[[''B'', ''E'', ''B''], [''A'', ''B'', ''E'']]'
$ws.Range("E8").Value = '[[''E'', ''A'', ''E''], [''D'', ''E'', ''A'']]'
$ws.Range("F8").Value = '[(46.272131, 52.17), (44.310045, 50.196303)]'
$ws.Range("G8").Value = '[(13.878853, 21.111869), (60.040123, 71.429511)]'
$ws.Range("H8").Value = 'spotify:track:5SUlhldQJtOhUr2GzH5RI7'
$ws.Range("I8").Value = ""

$ws.Range("A9").Value = 'isophonics_136'
$ws.Range("B9").Value = 'isophonics_288'
$ws.Range("C9").Value = 0.0782520325203252
$ws.Range("D9").Value = '[[''E/4'', ''D/5'', ''A''], [''E/3'', ''A'', ''D'']]'
$ws.Range("E9").Value = '[[''B'', ''A'', ''E''], [''B'', ''E'', ''A'']]'
$ws.Range("F9").Value = '[(32.311, 37.28), (95.673, 99.127)]'
$ws.Range("G9").Value = '[(0.268067, 8.504874), (30.424512, 35.962471)]'
$ws.Range("H9").Value = ""
$ws.Range("I9").Value = ""

$ws.Range("A10").Value = 'schubert-winterreise_65'
$ws.Range("B10").Value = 'isophonics_233'
$ws.Range("C10").Value = 0.1441176470588235
$ws.Range("D10").Value = '[[''D#:min'', ''A#:7'', ''D#:min'']]'
$ws.Range("E10").Value = '[[''D:min'', ''A:7'', ''D:min'']]'
$ws.Range("F10").Value = '[(9.1, 13.86)]'
$ws.Range("G10").Value = '[(36.711, 39.787)]'
$ws.Range("H10").Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'
$ws.Range("I10").Value = ""

$ws.Range("A11").Value = 'isophonics_234'
$ws.Range("B11").Value = 'schubert-winterreise_145'
$ws.Range("C11").Value = 0.4666666666666667
$ws.Range("D11").Value = '[[''Eb/5'', ''Ab/2'', ''Eb'']]'
$ws.Range("E11").Value = '[[''D:maj/A'', ''G:maj'', ''D:maj/A'']]'
$ws.Range("F11").Value = '[(75.439, 83.381)]'
$ws.Range("G11").Value = '[(143.58, 148.16)]'
$ws.Range("H11").Value = ""
$ws.Range("I11").Value = ""

$ws.Range("A12").Value = 'isophonics_69'
$ws.Range("B12").Value = 'isophonics_221'
$ws.Range("C12").Value = 0.2549019607843137
$ws.Range("D12").Value = '[[''E/5'', ''E/#4'', ''E/4''], [''E'', ''A'', ''E'']]'
$ws.Range("E12").Value = '[[''C'', ''C/7'', ''C/6''], [''C'', ''F/5'', ''C'']]'
$ws.Range("F12").Value = '[(0.440395, 4.837819), (25.015959, 36.625937)]'
$ws.Range("G12").Value = '[(69.901, 77.903), (5.121, 12.601)]'
$ws.Range("H12").Value = ""
$ws.Range("I12").Value = ""

$ws.Range("A13").Value = 'isophonics_273'
$ws.Range("B13").Value = 'isophonics_297'
$ws.Range("C13").Value = 0.08492822966507177
$ws.Range("D13").Value = '[[''C'', ''G'', ''C'']]'
$ws.Range("E13").Value = '[[''G'', ''D'', ''G'']]'
$ws.Range("F13").Value = '[(20.801, 27.165)]'
$ws.Range("G13").Value = '[(0.421247, 3.083177)]'
$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

$ws.Range("A14").Value = 'jaah_39'
$ws.Range("B14").Value = 'isophonics_279'
$ws.Range("C14").Value = 0.05900948366701791
$ws.Range("D14").Value = '[[''Db'', ''Db:min'', ''Ab'']]'
$ws.Range("E14").Value = '[[''F:maj'', ''F:min'', ''C:maj'']]'
$ws.Range("F14").Value = '[(129.12, 131.12)]'
$ws.Range("G14").Value = '[(46.097, 54.196)]'
$ws.Range("H14").Value = ""
$ws.Range("I14").Value = ""

$ws.Range("A15").Value = 'isophonics_45'
$ws.Range("B15").Value = 'schubert-winterreise_84'
$ws.Range("C15").Value = 0.1396103896103896
$ws.Range("D15").Value = '[[''F:maj'', ''F:7'', ''Bb'']]'
$ws.Range("E15").Value = '[[''D#:maj'', ''D#:7'', ''G#:maj'']]'
$ws.Range("F15").Value = '[(13.155, 20.379)]'
$ws.Range("G15").Value = '[(15.64, 24.28)]'
$ws.Range("H15").Value = ""
$ws.Range("I15").Value = ""

$ws.Range("A16").Value = 'isophonics_196'
$ws.Range("B16").Value = 'isophonics_170'
$ws.Range("C16").Value = 0.07668898356664638
$ws.Range("D16").Value = '[[''A'', ''D/5'', ''D'']]'
$ws.Range("E16").Value = '[[''A'', ''D'', ''D/7'']]'
$ws.Range("F16").Value = '[(35.463242, 38.330907)]'
$ws.Range("G16").Value = '[(10.560212, 12.866913)]'
$ws.Range("H16").Value = ""
$ws.Range("I16").Value = ""

$ws.Range("A17").Value = 'schubert-winterreise_177'
$ws.Range("B17").Value = 'schubert-winterreise_143'
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = '[[''A:min'', ''A:7'', ''D:maj/A'', ''D:min/A'', ''A:maj'', ''A:min'', ''E:maj'', ''A:min'', ''A:7'', ''D:min/A'', ''E:7/A'', ''A:min'', ''D:min/A'', ''A:min'', ''C:maj/G'', ''D:min/F'', ''C:maj/E'', ''G:7'', ''C:maj'', ''A:min'', ''A:7/G'', ''D:maj/F#'', ''G:(3,5,b7,b9)/F'', ''C:maj/E'', ''E:(3,5,b7,b9)/D'', ''F:(3,5)'', ''A:min/E'', ''E:7'', ''F:(3,5)'', ''A:min/E'', ''E:7'', ''A:min'', ''A:7'', ''D:maj'', ''D:min'', ''A:maj'', ''A:min/A'', ''E:maj'', ''A:min'']]'
$ws.Range("E17").Value = '[[''B:min'', ''B:7'', ''E:maj/B'', ''E:min/B'', ''B:maj'', ''B:min'', ''F#:maj'', ''B:min'', ''B:7'', ''E:min/B'', ''F#:7/B'', ''B:min'', ''E:min/B'', ''B:min'', ''D:maj/A'', ''E:min/G'', ''D:maj/F#'', ''A:7'', ''D:maj'', ''B:min'', ''B:7/A'', ''E:maj/G#'', ''A:(3,5,b7,b9)/G'', ''D:maj/F#'', ''F#:(3,5,b7,b9)/E'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''G:(3,5)'', ''B:min/F#'', ''F#:7'', ''B:min'', ''B:7'', ''E:maj'', ''E:min'', ''B:maj'', ''B:min/B'', ''F#:maj'', ''B:min'']]'
$ws.Range("F17").Value = '[(0.62, 107.38)]'
$ws.Range("G17").Value = '[(1.66, 97.0)]'
$ws.Range("H17").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
$ws.Range("I17").Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'
